# Se procesan de nuevo los datos con las nuevas dimensiones curadas
#
# Columns A (temporalidad), G (sector-descripcion), K (mes-nombre),
# N (dias-duracion-contrato) and P (sexo) move from "dimension" to
# "measure" semantics: row2 label, row3 kind ("dim"->"medida"), row4
# datatype ("skos:Concept"->"xsd:int"), and the row5 mapping-file
# reference is dropped (those columns no longer need a curated mapping
# workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("A", "G", "K", "N", "P")

foreach ($col in $columns) {
    $row2 = $ws.Range($col + "2")
    $oldLabel = $row2.Value2
    $row2.Value = $oldLabel -replace "^iaest-dimension:", "iaest-measure:"

    $ws.Range($col + "3").Value = "medida"
    $ws.Range($col + "4").Value = "xsd:int"
    $ws.Range($col + "5").ClearContents()
}
